# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型"
# sheets to match the latest scrape, as described in the commit:
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (column F holds 想去人数) ---
$wsExh = $wb.Worksheets.Item("展览")
$wsExh.Cells.Item(2, 6).Value  = 590   # 南宁·小蜜蜂动漫嘉年华2.0        588 -> 590
$wsExh.Cells.Item(5, 6).Value  = 479   # 南宁·0713国乙ONLY              473 -> 479
$wsExh.Cells.Item(7, 6).Value  = 2514  # 南宁·AB动漫游戏嘉年华          2511 -> 2514
$wsExh.Cells.Item(9, 6).Value  = 6716  # 南宁·良牙动漫夏季盛典          6700 -> 6716
$wsExh.Cells.Item(10, 6).Value = 179   # 南宁·火影忍者only              178 -> 179
$wsExh.Cells.Item(11, 6).Value = 428   # 南宁·蔚蓝档案only              427 -> 428

# --- Sheet "全部类型" (column F holds 想去人数) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Cells.Item(2, 6).Value  = 590   # 南宁·小蜜蜂动漫嘉年华2.0        588 -> 590
$wsAll.Cells.Item(5, 6).Value  = 479   # 南宁·0713国乙ONLY              473 -> 479
$wsAll.Cells.Item(9, 6).Value  = 2514  # 南宁·AB动漫游戏嘉年华          2511 -> 2514
$wsAll.Cells.Item(11, 6).Value = 6716  # 南宁·良牙动漫夏季盛典          6700 -> 6716
$wsAll.Cells.Item(12, 6).Value = 179   # 南宁·火影忍者only              178 -> 179
$wsAll.Cells.Item(13, 6).Value = 428   # 南宁·蔚蓝档案only              427 -> 428
